$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.208.30"
$ws.Range("E2").Value = "  +0.18%  "
$ws.Range("D3").Value = "3.517.76"
$ws.Range("E3").Value = "  +0.09%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "595.41"
$ws.Range("E5").Value = "  +0.29%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "172.60"
$ws.Range("E6").Value = "  +1.64%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("E8").Value = "  +2.37%  "
$ws.Range("E9").Value = "  +6.98%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.30"
$ws.Range("E10").Value = "  +0.46%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.435"
$ws.Range("E11").Value = "  -0.91%  "
$ws.Range("D12").Value = "4.133.94"
$ws.Range("E12").Value = "  +0.13%  "
$ws.Range("E13").Value = "  -0.21%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.68"
$ws.Range("E14").Value = "  +1.46%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000181"
$ws.Range("E15").Value = "  +0.58%  "
$ws.Range("D16").Value = "67.164.61"
$ws.Range("E16").Value = "  +0.25%  "
$ws.Range("D17").Value = "3.510.33"
$ws.Range("E17").Value = "  -0.49%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.34"
$ws.Range("E18").Value = "  +0.58%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "14.23"
$ws.Range("E19").Value = "  +1.07%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "396.21"
$ws.Range("E20").Value = "  +0.00%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.98"
$ws.Range("E21").Value = "  +0.63%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "73.40"
$ws.Range("E22").Value = "  -0.19%  "
$ws.Range("B23").Value = "Dai"
$ws.Range("C23").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("E23").Value = "  +0.05%  "
$ws.Range("B24").Value = "Polygon"
$ws.Range("C24").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.539"
$ws.Range("E24").Value = "  +1.56%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000122"
$ws.Range("E25").Value = "  -4.28%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.30"
$ws.Range("E26").Value = "  +2.00%  "
$ws.Range("E27").Value = "  -0.99%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.998"
$ws.Range("E28").Value = "  -0.21%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.28"
$ws.Range("E29").Value = "  -1.92%  "
$ws.Range("E30").Value = "  -0.85%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.08"
$ws.Range("E31").Value = "  +0.72%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "24.09"
$ws.Range("E32").Value = "  +1.95%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.40"
$ws.Range("E33").Value = "  +0.07%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.63"
$ws.Range("E34").Value = "  +2.43%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "163.77"
$ws.Range("E35").Value = "  +1.28%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.895"
$ws.Range("E36").Value = "  -1.26%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.91"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.89"
$ws.Range("E38").Value = "  +2.41%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.72"
$ws.Range("E39").Value = "  +0.76%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0746"
$ws.Range("E40").Value = "  -0.51%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "27.38"
$ws.Range("E41").Value = "  +1.54%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "26.38"
$ws.Range("E42").Value = "  +0.05%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.62"
$ws.Range("E43").Value = "  +2.31%  "
$ws.Range("D44").Value = "2.808.89"
$ws.Range("E44").Value = "  -0.72%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "42.94"
$ws.Range("E45").Value = "  -1.34%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0309"
$ws.Range("E46").Value = "  -2.14%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "341.84"
$ws.Range("E47").Value = "  -3.27%  "
$ws.Range("E48").Value = "  +0.66%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "33.94"
$ws.Range("E49").Value = "  +1.96%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.53"
$ws.Range("E50").Value = "  +0.57%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.850"
$ws.Range("E51").Value = "  -0.25%  "
